$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (B1:N1) - replace "/Category/Subcategory" style with "Category-Subcategory" (underscored)
$ws.Range("B1").Value = "Government-Cadw"
$ws.Range("C1").Value = "Government-Local_Authority"
$ws.Range("D1").Value = "Government-National"
$ws.Range("E1").Value = "Government-Other"
$ws.Range("F1").Value = "Independent-English_Heritage"
$ws.Range("G1").Value = "Independent-Historic_Environment_Scotland"
$ws.Range("H1").Value = "Independent-National_Trust"
$ws.Range("I1").Value = "Independent-National_Trust_for_Scotland"
$ws.Range("J1").Value = "Independent-Not_for_profit"
$ws.Range("K1").Value = "Independent-Private"
$ws.Range("L1").Value = "Independent-Unknown"
$ws.Range("M1").Value = "University"
$ws.Range("N1").Value = "Unknown"

# Update numeric data values for row 2 (government)
$ws.Range("C2").Value = 922
$ws.Range("D2").Value = 82
$ws.Range("E2").Value = 10

# Update numeric data values for row 3 (independent)
$ws.Range("F3").Value = 53
$ws.Range("G3").Value = 21
$ws.Range("I3").Value = 27
$ws.Range("J3").Value = 1734
$ws.Range("K3").Value = 751
$ws.Range("L3").Value = 221

# Update numeric data value for row 4 (university)
$ws.Range("M4").Value = 110

# Update numeric data value for row 5 (unknown_gov)
$ws.Range("N5").Value = 110
